$d = $word.ActiveDocument

# The document's last paragraph is the final list item of "Tópico 2"
# ("– Upload/Verificar_Status_de_Chamado", style PargrafodaLista,
# ilvl 1 / numId 2). Append a new sibling list item, "Delete", right after
# it (and right before the sectPr), inheriting the same list formatting.
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "Delete"
